{"js": "// Update the 25 \"NNN\u00d7N=\" expression cells in the multiplication-practice\n// table from their old values to the newly generated ones.\nconst replacements = [\n  [\"175\u00d79=\", \"557\u00d77=\"],\n  [\"184\u00d75=\", \"183\u00d77=\"],\n  [\"523\u00d79=\", \"505\u00d79=\"],\n  [\"798\u00d77=\", \"776\u00d79=\"],\n  [\"213\u00d78=\", \"512\u00d73=\"],\n  [\"306\u00d74=\", \"750\u00d75=\"],\n  [\"582\u00d72=\", \"989\u00d75=\"],\n  [\"684\u00d72=\", \"129\u00d79=\"],\n  [\"154\u00d78=\", \"391\u00d77=\"],\n  [\"510\u00d79=\", \"132\u00d79=\"],\n  [\"262\u00d73=\", \"351\u00d73=\"],\n  [\"972\u00d78=\", \"638\u00d74=\"],\n  [\"823\u00d75=\", \"946\u00d73=\"],\n  [\"122\u00d73=\", \"444\u00d78=\"],\n  [\"526\u00d73=\", \"150\u00d79=\"],\n  [\"853\u00d73=\", \"424\u00d73=\"],\n  [\"829\u00d75=\", \"316\u00d76=\"],\n  [\"127\u00d72=\", \"361\u00d78=\"],\n  [\"796\u00d76=\", \"257\u00d74=\"],\n  [\"145\u00d79=\", \"924\u00d72=\"],\n  [\"576\u00d73=\", \"381\u00d79=\"],\n  [\"835\u00d76=\", \"254\u00d73=\"],\n  [\"361\u00d72=\", \"807\u00d79=\"],\n  [\"541\u00d77=\", \"914\u00d72=\"],\n  [\"603\u00d75=\", \"152\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 \"NNN\u00d7N=\" expression cells in the multiplication-practice\n# table from their old values to the newly generated ones.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"175\u00d79=\", \"557\u00d77=\"),\n    @(\"184\u00d75=\", \"183\u00d77=\"),\n    @(\"523\u00d79=\", \"505\u00d79=\"),\n    @(\"798\u00d77=\", \"776\u00d79=\"),\n    @(\"213\u00d78=\", \"512\u00d73=\"),\n    @(\"306\u00d74=\", \"750\u00d75=\"),\n    @(\"582\u00d72=\", \"989\u00d75=\"),\n    @(\"684\u00d72=\", \"129\u00d79=\"),\n    @(\"154\u00d78=\", \"391\u00d77=\"),\n    @(\"510\u00d79=\", \"132\u00d79=\"),\n    @(\"262\u00d73=\", \"351\u00d73=\"),\n    @(\"972\u00d78=\", \"638\u00d74=\"),\n    @(\"823\u00d75=\", \"946\u00d73=\"),\n    @(\"122\u00d73=\", \"444\u00d78=\"),\n    @(\"526\u00d73=\", \"150\u00d79=\"),\n    @(\"853\u00d73=\", \"424\u00d73=\"),\n    @(\"829\u00d75=\", \"316\u00d76=\"),\n    @(\"127\u00d72=\", \"361\u00d78=\"),\n    @(\"796\u00d76=\", \"257\u00d74=\"),\n    @(\"145\u00d79=\", \"924\u00d72=\"),\n    @(\"576\u00d73=\", \"381\u00d79=\"),\n    @(\"835\u00d76=\", \"254\u00d73=\"),\n    @(\"361\u00d72=\", \"807\u00d79=\"),\n    @(\"541\u00d77=\", \"914\u00d72=\"),\n    @(\"603\u00d75=\", \"152\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n}\n"}
